$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename header cells ---
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after "Monthly Trend" ---
$new = $wb.Worksheets.Add($null, $ws2)
$new.Name = "PO Forecast"

# Copy header style (bold/border/centered) from an existing header row
$ws1.Range("A1:B1").Copy()
$new.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-number-format style from column A data cells
$ws1.Range("A2").Copy()
$new.Range("A2:A27").PasteSpecial(-4122)

$new.Range("A1").Value = "ds"
$new.Range("B1").Value = "PO_Forecast"
$new.Range("C1").Value = "yhat_lower"
$new.Range("D1").Value = "yhat_upper"

$data = @(
    @(45424.99999999999, 34, -1.101575294931065, 69.34100575817709),
    @(45445.99999999999, 36, 2.471164277202467, 72.57888785513055),
    @(45459.99999999999, 37, 3.015830122124771, 73.076973424157),
    @(45466.99999999999, 38, 0.5774569914005367, 73.62299833926775),
    @(45473.99999999999, 39, 0.1603969116946179, 73.40539477549659),
    @(45480.99999999999, 39, 4.954054558003199, 73.82694142651799),
    @(45487.99999999999, 40, 6.829606299366305, 74.73621383781339),
    @(45494.99999999999, 41, 4.007788115915539, 77.21660776781495),
    @(45508.99999999999, 42, 6.249013504470854, 76.13671947397582),
    @(45515.99999999999, 42, 7.647751644207788, 77.56787030197825),
    @(45522.99999999999, 43, 6.10500301935176, 78.12710506093343),
    @(45529.99999999999, 44, 8.78040608265883, 78.10006353012172),
    @(45550.99999999999, 46, 9.089537848473316, 80.34900536360811),
    @(45557.99999999999, 46, 10.20145134424669, 83.12248449084471),
    @(45571.99999999999, 48, 14.85512551984903, 83.74859436356336),
    @(45578.99999999999, 48, 14.71898837644669, 83.90512770817553),
    @(45585.99999999999, 49, 14.77125500309396, 86.89641574074464),
    @(45592.99999999999, 49, 12.35714399086269, 84.43123787826966),
    @(45599.99999999999, 50, 11.98972951371482, 85.26162234284429),
    @(45606.99999999999, 51, 14.22618215339896, 86.37653071384349),
    @(45613.99999999999, 51, 17.48404096415989, 87.37332788197553),
    @(45620.99999999999, 52, 15.89319715730741, 88.39123981236874),
    @(45627.99999999999, 53, 18.70157432984298, 88.50130600558266),
    @(45634.99999999999, 53, 15.7195759756893, 90.03714333571145),
    @(45641.99999999999, 54, 20.50933734863232, 90.55330657005321),
    @(45648.99999999999, 54, 20.0998940673528, 88.6711998738605)
)

$r = 2
foreach ($row in $data) {
    $new.Cells.Item($r, 1).Value = $row[0]
    $new.Cells.Item($r, 2).Value = $row[1]
    $new.Cells.Item($r, 3).Value = $row[2]
    $new.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Output "done"
